$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ebi3"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 2.878032333333334
$ws.Range("H2").Value = 8.634097
$ws.Range("I2").Value = 0.3279446910817746
$ws.Range("J2").Value = 0.3279446910817746
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 36.71344366666667
$ws.Range("N2").Value = 110.140331
$ws.Range("O2").Value = 0.2081992981130139
$ws.Range("P2").Value = 0.2081992981130138
$ws.Range("Q2").Value = 105.6624779406786
$ws.Range("R2").Value = 950.9623014661071
$ws.Range("S2").Value = 0.06827785450311462
$ws.Range("T2").Value = 0.0682778545031146

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ebi3"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 2.878032333333334
$ws.Range("H3").Value = 8.634097
$ws.Range("I3").Value = 0.3279446910817746
$ws.Range("J3").Value = 0.3279446910817746
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 88.87708033333332
$ws.Range("N3").Value = 266.631241
$ws.Range("O3").Value = 0.5040155293450301
$ws.Range("P3").Value = 0.50401552934503
$ws.Range("Q3").Value = 255.7911108915974
$ws.Range("R3").Value = 2302.119998024377
$ws.Range("S3").Value = 0.165289217071473
$ws.Range("T3").Value = 0.165289217071473

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ebi3"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 2.878032333333334
$ws.Range("H4").Value = 8.634097
$ws.Range("I4").Value = 0.3279446910817746
$ws.Range("J4").Value = 0.3279446910817746
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 14.22727166666667
$ws.Range("N4").Value = 42.681815
$ws.Range("O4").Value = 0.08068183420648613
$ws.Range("P4").Value = 0.08068183420648613
$ws.Range("Q4").Value = 40.94654787178389
$ws.Range("R4").Value = 368.518930846055
$ws.Range("S4").Value = 0.02645917919475705
$ws.Range("T4").Value = 0.02645917919475705

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ebi3"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 2.878032333333334
$ws.Range("H5").Value = 8.634097
$ws.Range("I5").Value = 0.3279446910817746
$ws.Range("J5").Value = 0.3279446910817746
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 36.52018433333333
$ws.Range("N5").Value = 109.560553
$ws.Range("O5").Value = 0.20710333833547
$ws.Range("P5").Value = 0.2071033383354699
$ws.Range("Q5").Value = 105.1062713306268
$ws.Range("R5").Value = 945.956441975641
$ws.Range("S5").Value = 0.06791844031242994
$ws.Range("T5").Value = 0.06791844031242993

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ebi3"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 5.897936333333334
$ws.Range("H6").Value = 17.693809
$ws.Range("I6").Value = 0.6720553089182254
$ws.Range("J6").Value = 0.6720553089182254
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 36.71344366666667
$ws.Range("N6").Value = 110.140331
$ws.Range("O6").Value = 0.2081992981130139
$ws.Range("P6").Value = 0.2081992981130138
$ws.Range("Q6").Value = 216.5335533234199
$ws.Range("R6").Value = 1948.801979910779
$ws.Range("S6").Value = 0.1399214436098992
$ws.Range("T6").Value = 0.1399214436098992

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ebi3"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 5.897936333333334
$ws.Range("H7").Value = 17.693809
$ws.Range("I7").Value = 0.6720553089182254
$ws.Range("J7").Value = 0.6720553089182254
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 88.87708033333332
$ws.Range("N7").Value = 266.631241
$ws.Range("O7").Value = 0.5040155293450301
$ws.Range("P7").Value = 0.50401552934503
$ws.Range("Q7").Value = 524.1913612985521
$ws.Range("R7").Value = 4717.72225168697
$ws.Range("S7").Value = 0.3387263122735572
$ws.Range("T7").Value = 0.338726312273557

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ebi3"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 5.897936333333334
$ws.Range("H8").Value = 17.693809
$ws.Range("I8").Value = 0.6720553089182254
$ws.Range("J8").Value = 0.6720553089182254
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 14.22727166666667
$ws.Range("N8").Value = 42.681815
$ws.Range("O8").Value = 0.08068183420648613
$ws.Range("P8").Value = 0.08068183420648613
$ws.Range("Q8").Value = 83.91154248703722
$ws.Range("R8").Value = 755.2038823833351
$ws.Range("S8").Value = 0.05422265501172908
$ws.Range("T8").Value = 0.05422265501172908

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ebi3"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 5.897936333333334
$ws.Range("H9").Value = 17.693809
$ws.Range("I9").Value = 0.6720553089182254
$ws.Range("J9").Value = 0.6720553089182254
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 36.52018433333333
$ws.Range("N9").Value = 109.560553
$ws.Range("O9").Value = 0.20710333833547
$ws.Range("P9").Value = 0.2071033383354699
$ws.Range("Q9").Value = 215.3937220795975
$ws.Range("R9").Value = 1938.543498716377
$ws.Range("S9").Value = 0.13918489802304
$ws.Range("T9").Value = 0.13918489802304
